# Refresh currentAveragePrice / Leve cost & profit columns (H-N) across all crafting-leve sheets
# to match the latest market-board snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 4067613.2
$ws.Range("I70").Value = 6098670
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 18296010
$ws.Range("L70").Value = 16500
$ws.Range("M70").Value = -18295740
$ws.Range("N70").Value = -17040
# Row 73
$ws.Range("H73").Value = 4067613.2
$ws.Range("I73").Value = 6098670
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 18296010
$ws.Range("L73").Value = 16500
$ws.Range("M73").Value = -18295074
$ws.Range("N73").Value = -18372
# Row 94
$ws.Range("H94").Value = 3607.9
$ws.Range("I94").Value = 3607.9
$ws.Range("K94").Value = 3607.9
$ws.Range("M94").Value = -3156.9
# Row 100
$ws.Range("H100").Value = 4588.7915
$ws.Range("I100").Value = 3040.1177
$ws.Range("K100").Value = 3040.1177
$ws.Range("M100").Value = -2499.1177
# Row 103
$ws.Range("I103").Value = 1499
$ws.Range("J103").Value = 50001876
$ws.Range("K103").Value = 4497
$ws.Range("L103").Value = 150005628
$ws.Range("M103").Value = -3911
$ws.Range("N103").Value = -150006800
# Row 106
$ws.Range("H106").Value = 7794.1816
$ws.Range("I106").Value = 7794.1816
$ws.Range("K106").Value = 7794.1816
$ws.Range("M106").Value = -7163.1816
# Row 138
$ws.Range("H138").Value = 5192.1836
$ws.Range("I138").Value = 3096.75
$ws.Range("J138").Value = 6208.1514
$ws.Range("K138").Value = 9290.25
$ws.Range("L138").Value = 18624.4542
$ws.Range("M138").Value = -4150.25
$ws.Range("N138").Value = -28904.4542

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10830.951
$ws.Range("I32").Value = 10587.509
$ws.Range("J32").Value = 14300
$ws.Range("K32").Value = 10587.509
$ws.Range("L32").Value = 14300
$ws.Range("M32").Value = -10300.509
$ws.Range("N32").Value = -14874
# Row 88
$ws.Range("H88").Value = 2933.6667
$ws.Range("J88").Value = 2890.5
$ws.Range("L88").Value = 2890.5
$ws.Range("N88").Value = -3702.5
# Row 91
$ws.Range("H91").Value = 2933.6667
$ws.Range("J91").Value = 2890.5
$ws.Range("L91").Value = 2890.5
$ws.Range("N91").Value = -5698.5
# Row 132
$ws.Range("H132").Value = 5005354.5
$ws.Range("I132").Value = 4735.467
$ws.Range("K132").Value = 14206.401
$ws.Range("M132").Value = -11676.401

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1125.7858
$ws.Range("J86").Value = 2220
$ws.Range("L86").Value = 2220
$ws.Range("N86").Value = -4466
# Row 89
$ws.Range("H89").Value = 1125.7858
$ws.Range("J89").Value = 2220
$ws.Range("L89").Value = 11100
$ws.Range("N89").Value = -22332
# Row 105
$ws.Range("H105").Value = 381014.28
$ws.Range("I105").Value = 515834.88
$ws.Range("J105").Value = 6512.6665
$ws.Range("K105").Value = 515834.88
$ws.Range("L105").Value = 6512.6665
$ws.Range("M105").Value = -514087.88
$ws.Range("N105").Value = -10006.6665
# Row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").ClearContents()
$ws.Range("N108").Value = 0

$ws = $wb.Worksheets.Item("CRP")
# Row 107
$ws.Range("H107").Value = 1505.65
$ws.Range("I107").Value = 482
$ws.Range("J107").Value = 3041.125
$ws.Range("K107").Value = 482
$ws.Range("L107").Value = 3041.125
$ws.Range("M107").Value = 1438
$ws.Range("N107").Value = -6881.125
# Row 134
$ws.Range("H134").Value = 2517.087
$ws.Range("I134").Value = 2219.65
$ws.Range("K134").Value = 6658.950000000001
$ws.Range("M134").Value = -4123.950000000001
# Row 141
$ws.Range("H141").Value = 522925.34
$ws.Range("J141").Value = 592510.4
$ws.Range("L141").Value = 592510.4
$ws.Range("N141").Value = -602870.4

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1687.2858
$ws.Range("I113").Value = 1412.8334
$ws.Range("J113").Value = 2053.2222
$ws.Range("K113").Value = 4238.5002
$ws.Range("L113").Value = 6159.6666
$ws.Range("M113").Value = -2068.5002
$ws.Range("N113").Value = -10499.6666

$ws = $wb.Worksheets.Item("GSM")
# Row 136
$ws.Range("H136").Value = 9482.277
$ws.Range("J136").Value = 9482.277
$ws.Range("L136").Value = 28446.831
$ws.Range("N136").Value = -33546.831

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 9997
$ws.Range("J2").Value = 9997
$ws.Range("L2").Value = 9997
$ws.Range("N2").Value = -10221
# Row 100
$ws.Range("H100").Value = 16670918
$ws.Range("I100").Value = 4274
$ws.Range("J100").Value = 41670884
$ws.Range("K100").Value = 4274
$ws.Range("L100").Value = 41670884
$ws.Range("M100").Value = -3733
$ws.Range("N100").Value = -41671966

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8035.174
$ws.Range("I62").Value = 5040.4165
$ws.Range("J62").Value = 11302.182
$ws.Range("K62").Value = 5040.4165
$ws.Range("L62").Value = 11302.182
$ws.Range("M62").Value = -4416.4165
$ws.Range("N62").Value = -12550.182
# Row 65
$ws.Range("H65").Value = 8035.174
$ws.Range("I65").Value = 5040.4165
$ws.Range("J65").Value = 11302.182
$ws.Range("K65").Value = 25202.0825
$ws.Range("L65").Value = 56510.91
$ws.Range("M65").Value = -22082.0825
$ws.Range("N65").Value = -62750.91
# Row 81
$ws.Range("H81").Value = 1427.3334
$ws.Range("I81").Value = 1303.1666
$ws.Range("J81").Value = 1924
$ws.Range("K81").Value = 2606.3332
$ws.Range("L81").Value = 3848
$ws.Range("M81").Value = -1545.3332
$ws.Range("N81").Value = -5970
# Row 84
$ws.Range("H84").Value = 1427.3334
$ws.Range("I84").Value = 1303.1666
$ws.Range("J84").Value = 1924
$ws.Range("K84").Value = 13031.666
$ws.Range("L84").Value = 19240
$ws.Range("M84").Value = -7727.666000000001
$ws.Range("N84").Value = -29848
# Row 107
$ws.Range("H107").Value = 2788.98
$ws.Range("I107").Value = 1442.4865
$ws.Range("J107").Value = 6621.3076
$ws.Range("K107").Value = 4327.4595
$ws.Range("L107").Value = 19863.9228
$ws.Range("M107").Value = -2407.4595
$ws.Range("N107").Value = -23703.9228
# Row 122
$ws.Range("H122").Value = 3108.389
$ws.Range("I122").Value = 2919.6924
$ws.Range("J122").Value = 3599
$ws.Range("K122").Value = 8759.0772
$ws.Range("L122").Value = 10797
$ws.Range("M122").Value = -6309.0772
$ws.Range("N122").Value = -15697
# Row 141
$ws.Range("H141").Value = 95998.664
$ws.Range("J141").Value = 95998.664
$ws.Range("L141").Value = 95998.664
$ws.Range("N141").Value = -106358.664
